# AHC-20-001 DP Switch Barrier Block Interface BOM - rev 1 final fab files
# Test points J5/J6 were replaced by TP1-TP4 on the final fab revision.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Item 4's "Part Reference" cell (C5) changes from "J5 J6" to "TP1 TP2 TP3 TP4"
$ws.Range("C5").Value = "TP1 TP2 TP3 TP4"

# Column C widens (still autosized / "best fit") to accommodate the longer text
$ws.Columns.Item(3).ColumnWidth = 14.1

# Reset the lingering UI selection back to the sheet's home cell
$ws.Range("A1").Select()
